$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be treated as text so purely-numeric-looking
    # strings (e.g. "318.48") are not silently coerced into floating
    # point numbers by Excel's automatic type detection.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "44.084.00"
$ws.Range("E2").Value = "  +2.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.255.26"
$ws.Range("E3").Value = "  +1.65%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.23%  "

# Row 5 - BNB
Set-TextValue "D5" "318.48"
$ws.Range("E5").Value = "  -0.09%  "

# Row 6 - Solana
Set-TextValue "D6" "101.22"
$ws.Range("E6").Value = "  +3.13%  "

# Row 7 - XRP
Set-TextValue "D7" "0.579"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.15%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.552"
$ws.Range("E9").Value = "  -2.17%  "

# Row 10 - Avalanche
Set-TextValue "D10" "37.09"
$ws.Range("E10").Value = "  +1.64%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0837"
$ws.Range("E11").Value = "  +1.47%  "

# Row 12 - Polkadot
Set-TextValue "D12" "7.58"
$ws.Range("E12").Value = "  -0.23%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.32%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.602.65"
$ws.Range("E14").Value = "  +2.06%  "

# Row 15 - Chainlink
Set-TextValue "D15" "14.53"
$ws.Range("E15").Value = "  +1.75%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.859"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.261.28"
$ws.Range("E17").Value = "  +1.82%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.976.49"
$ws.Range("E18").Value = "  +2.66%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "13.36"
$ws.Range("E19").Value = "  -2.32%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  +2.59%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.48"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22 - Litecoin
Set-TextValue "D22" "65.73"
$ws.Range("E22").Value = "  +1.03%  "

# Row 23 - PancakeSwap
$ws.Range("E23").Value = "  -3.11%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "235.14"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25 - ImmutableX
Set-TextValue "D25" "2.08"
$ws.Range("E25").Value = "  -4.69%  "

# Row 27 - Cosmos
Set-TextValue "D27" "10.62"
$ws.Range("E27").Value = "  +6.54%  "

# Row 28 - InjectiveProtocol
Set-TextValue "D28" "38.55"
$ws.Range("E28").Value = "  +6.65%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -1.62%  "

# Row 30 - Filecoin
Set-TextValue "D30" "6.19"
$ws.Range("E30").Value = "  -2.37%  "

# Row 31 - Monero
Set-TextValue "D31" "161.82"
$ws.Range("E31").Value = "  +4.42%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "20.20"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0848"
$ws.Range("E33").Value = "  -1.43%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  +1.19%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.97"
$ws.Range("E35").Value = "  +6.02%  "

# Row 36 - Kaspa
Set-TextValue "D36" "0.113"
$ws.Range("E36").Value = "  +7.81%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -5.34%  "

# Row 38 - Stellar
$ws.Range("E38").Value = "  -0.84%  "

# Row 39 - Celestia
Set-TextValue "D39" "16.65"
$ws.Range("E39").Value = "  +19.54%  "

# Row 40 - was NEARProtocol, now RenderToken (rows 40/41 content swapped)
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D40" "4.19"
$ws.Range("E40").Value = "  -3.81%  "

# Row 41 - was RenderToken, now NEARProtocol
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D41" "3.66"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - VeChain
Set-TextValue "D42" "0.0316"
$ws.Range("E42").Value = "  -0.44%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.17%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.776.77"
$ws.Range("E44").Value = "  +2.63%  "

# Row 45 - Algorand
$ws.Range("E45").Value = "  -1.73%  "

# Row 46 - THORChain
Set-TextValue "D46" "5.21"
$ws.Range("E46").Value = "  -0.29%  "

# Row 47 - BitcoinSV
Set-TextValue "D47" "81.64"
$ws.Range("E47").Value = "  -2.39%  "

# Row 48 - ordi
Set-TextValue "D48" "74.43"
$ws.Range("E48").Value = "  +2.26%  "

# Row 49 - Aave
Set-TextValue "D49" "104.42"
$ws.Range("E49").Value = "  +1.92%  "

# Row 50 - was Stacks, now MultiversX (rows 50/51 content swapped)
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D50" "58.10"
$ws.Range("E50").Value = "  +1.34%  "

# Row 51 - was MultiversX, now Stacks
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D51" "1.67"
$ws.Range("E51").Value = "  +6.18%  "
